# Updating all workbooks for findings.
# For every worksheet: insert a new "auditee_uei" column at D (shifting the
# existing D:U columns to E:V), populate it with the auditee UEI, and
# convert the is_passthrough_award / is_repeat_finding columns (now at M
# and U respectively) from YES/NO text to native boolean TRUE/FALSE values.

$wb = $excel.ActiveWorkbook
$uei = "CXQZVRZCCF41"

foreach ($ws in $wb.Worksheets) {

    # 1) Insert the new column before column D; this shifts every existing
    #    column D..U one slot to the right (to E..V) including widths.
    $ws.Columns("D:D").Insert()

    # Give the new column the same width as the source file (~16.8 chars).
    $ws.Columns("D:D").ColumnWidth = 15.95

    # 2) Header + data values for the new auditee_uei column.
    $ws.Range("D1").Value = "auditee_uei"

    $lastRow = $ws.UsedRange.Rows.Count

    for ($r = 2; $r -le $lastRow; $r++) {
        $ws.Cells.Item($r, 4).Value = $uei

        # is_passthrough_award now lives in column M (13)
        $passCell = $ws.Cells.Item($r, 13)
        if ($passCell.Text -eq "YES") {
            $passCell.Value = $true
        } else {
            $passCell.Value = $false
        }
        $passCell.Style = "Normal"

        # is_repeat_finding now lives in column U (21)
        $repeatCell = $ws.Cells.Item($r, 21)
        if ($repeatCell.Text -eq "YES") {
            $repeatCell.Value = $true
        } else {
            $repeatCell.Value = $false
        }
        $repeatCell.Style = "Normal"
    }
}
